# Auto-generated Excel COM-interop script
# Applies numeric value updates to the H:N price/profit columns
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the commit diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 987.2308
$ws.Range("I4").Value = 985
$ws.Range("K4").Value = 985
$ws.Range("M4").Value = -871
$ws.Range("H18").Value = 9856.714
$ws.Range("I18").Value = 11332.833
$ws.Range("K18").Value = 11332.833
$ws.Range("M18").Value = -11048.833
$ws.Range("H43").Value = 8672.23
$ws.Range("J43").Value = 9328.25
$ws.Range("L43").Value = 9328.25
$ws.Range("N43").Value = -9466.25
$ws.Range("H100").Value = 61230.8
$ws.Range("I100").Value = 78355.08
$ws.Range("J100").Value = 42679.5
$ws.Range("K100").Value = 78355.08
$ws.Range("L100").Value = 42679.5
$ws.Range("M100").Value = -77814.08
$ws.Range("N100").Value = -43761.5
$ws.Range("H111").Value = 1154.4
$ws.Range("I111").Value = 944
$ws.Range("J111").Value = 1470
$ws.Range("K111").Value = 2832
$ws.Range("L111").Value = 4410
$ws.Range("M111").Value = 235
$ws.Range("N111").Value = -10544
$ws.Range("H116").Value = 6175723
$ws.Range("I116").Value = 9261549
$ws.Range("J116").Value = 4070
$ws.Range("K116").Value = 9261549
$ws.Range("L116").Value = 4070
$ws.Range("M116").Value = -9258107
$ws.Range("N116").Value = -10954
$ws.Range("H132").Value = 1788632.5
$ws.Range("I132").Value = 2928.157
$ws.Range("K132").Value = 8784.471000000001
$ws.Range("M132").Value = -6254.471000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7045.1665
$ws.Range("I45").Value = 6191
$ws.Range("K45").Value = 6191
$ws.Range("M45").Value = -5814
$ws.Range("H61").Value = 7758.5
$ws.Range("I61").Value = 8157.9585
$ws.Range("K61").Value = 8157.9585
$ws.Range("M61").Value = -7945.9585
$ws.Range("H74").Value = 1660.9656
$ws.Range("I74").Value = 939.35297
$ws.Range("K74").Value = 939.35297
$ws.Range("M74").Value = -65.35297000000003
$ws.Range("H77").Value = 1660.9656
$ws.Range("I77").Value = 939.35297
$ws.Range("K77").Value = 4696.76485
$ws.Range("M77").Value = -328.7648500000005
$ws.Range("H97").Value = 6901125.5
$ws.Range("I97").Value = 5707.381
$ws.Range("J97").Value = 25001598
$ws.Range("K97").Value = 5707.381
$ws.Range("L97").Value = 25001598
$ws.Range("M97").Value = -5211.381
$ws.Range("N97").Value = -25002590
$ws.Range("H110").Value = 2275.389
$ws.Range("I110").Value = 1674.5714
$ws.Range("K110").Value = 1674.5714
$ws.Range("M110").Value = 370.4286
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 40000
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 40000
$ws.Range("M114").ClearContents()
$ws.Range("N114").Value = -48678
$ws.Range("H136").Value = 7758.5
$ws.Range("I136").Value = 8157.9585
$ws.Range("K136").Value = 24473.8755
$ws.Range("M136").Value = -21923.8755

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 14718.333
$ws.Range("I82").Value = 5462
$ws.Range("J82").Value = 61000
$ws.Range("K82").Value = 5462
$ws.Range("L82").Value = 61000
$ws.Range("M82").Value = -5079
$ws.Range("N82").Value = -61766
$ws.Range("H85").Value = 14718.333
$ws.Range("I85").Value = 5462
$ws.Range("J85").Value = 61000
$ws.Range("K85").Value = 5462
$ws.Range("L85").Value = 61000
$ws.Range("M85").Value = -4136
$ws.Range("N85").Value = -63652
$ws.Range("H107").Value = 2580.8333
$ws.Range("J107").Value = 2866.3333
$ws.Range("L107").Value = 2866.3333
$ws.Range("N107").Value = -6706.3333
$ws.Range("H134").Value = 5851.2188
$ws.Range("I134").Value = 6162.0435
$ws.Range("K134").Value = 18486.1305
$ws.Range("M134").Value = -15951.1305

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 84880.836
$ws.Range("I16").Value = 1336.8334
$ws.Range("J16").Value = 168424.83
$ws.Range("K16").Value = 1336.8334
$ws.Range("L16").Value = 168424.83
$ws.Range("M16").Value = -1049.8334
$ws.Range("N16").Value = -168998.83
$ws.Range("H31").Value = 11163.4
$ws.Range("I31").Value = 14359.333
$ws.Range("K31").Value = 14359.333
$ws.Range("M31").Value = -14064.333
$ws.Range("H34").Value = 11163.4
$ws.Range("I34").Value = 14359.333
$ws.Range("K34").Value = 14359.333
$ws.Range("M34").Value = -14157.333
$ws.Range("H113").Value = 84880.836
$ws.Range("I113").Value = 1336.8334
$ws.Range("J113").Value = 168424.83
$ws.Range("K113").Value = 1336.8334
$ws.Range("L113").Value = 168424.83
$ws.Range("M113").Value = 833.1666
$ws.Range("N113").Value = -172764.83
$ws.Range("H120").Value = 61670
$ws.Range("J120").Value = 62219.6
$ws.Range("L120").Value = 62219.6
$ws.Range("N120").Value = -69477.60000000001
$ws.Range("H132").Value = 1518.4286
$ws.Range("I132").Value = 1400.9131
$ws.Range("K132").Value = 4202.7393
$ws.Range("M132").Value = -1672.7393
$ws.Range("H134").Value = 4668.2383
$ws.Range("I134").Value = 3638.75
$ws.Range("K134").Value = 10916.25
$ws.Range("M134").Value = -8381.25
$ws.Range("H141").Value = 194052.94
$ws.Range("J141").Value = 197769.11
$ws.Range("L141").Value = 197769.11
$ws.Range("N141").Value = -208129.11

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 500
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 500
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 1500
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -3496
$ws.Range("H78").Value = 500
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 500
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 4500
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -14484
$ws.Range("H137").Value = 4229.514
$ws.Range("I137").Value = 1851.625
$ws.Range("J137").Value = 9417.637000000001
$ws.Range("K137").Value = 5554.875
$ws.Range("L137").Value = 28252.911
$ws.Range("M137").Value = -454.875
$ws.Range("N137").Value = -38452.911

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 85000
$ws.Range("J108").Value = 85000
$ws.Range("L108").Value = 85000
$ws.Range("N108").Value = -92680
$ws.Range("H132").Value = 1896.4546
$ws.Range("I132").Value = 1914.025
$ws.Range("K132").Value = 5742.075000000001
$ws.Range("M132").Value = -3212.075000000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 736.63635
$ws.Range("I22").Value = 678.1111
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 678.1111
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -383.1111
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 736.63635
$ws.Range("I27").Value = 678.1111
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 678.1111
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -571.1111
$ws.Range("N27").Value = -1214
$ws.Range("H61").Value = 3692.963
$ws.Range("I61").Value = 1557.7142
$ws.Range("K61").Value = 1557.7142
$ws.Range("M61").Value = -1355.7142
$ws.Range("H68").Value = 9070.6
$ws.Range("I68").Value = 4000
$ws.Range("J68").Value = 10338.25
$ws.Range("K68").Value = 4000
$ws.Range("L68").Value = 10338.25
$ws.Range("M68").Value = -3251
$ws.Range("N68").Value = -11836.25
$ws.Range("H71").Value = 9070.6
$ws.Range("I71").Value = 4000
$ws.Range("J71").Value = 10338.25
$ws.Range("K71").Value = 20000
$ws.Range("L71").Value = 51691.25
$ws.Range("M71").Value = -16256
$ws.Range("N71").Value = -59179.25
$ws.Range("H74").Value = 32499
$ws.Range("I74").Value = 32499
$ws.Range("K74").Value = 32499
$ws.Range("M74").Value = -31501
$ws.Range("H77").Value = 32499
$ws.Range("I77").Value = 32499
$ws.Range("K77").Value = 97497
$ws.Range("M77").Value = -92505
$ws.Range("H113").Value = 3692.963
$ws.Range("I113").Value = 1557.7142
$ws.Range("K113").Value = 1557.7142
$ws.Range("M113").Value = 612.2858000000001
$ws.Range("H136").Value = 5253.6855
$ws.Range("I136").Value = 3769.5715
$ws.Range("K136").Value = 11308.7145
$ws.Range("M136").Value = -8758.7145

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10557.435
$ws.Range("I132").Value = 12105.151
$ws.Range("J132").Value = 6628.615
$ws.Range("K132").Value = 36315.453
$ws.Range("L132").Value = 19885.845
$ws.Range("M132").Value = -33785.453
$ws.Range("N132").Value = -24945.845

Write-Output "Applied all Siren_Profits value updates."